# Generate Report for Handback
# - row 3 (the "65f37892..." file) transitions from "Ready for handoff"
#   to "Handback transform failed" across Overview / zh-cn / de-de sheets
# - the zh-cn and de-de sheets get an "Error Detail" message in column P row 3
# - column P ("Error Detail") is widened to fit the new text

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handback transform failed"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

## Excel stores column widths internally as a pixel-rounded quantity derived
## from the "number of characters" ColumnWidth value; the source workbook's
## column 16 needs to land on an OOXML width of exactly 40. Empirically (via
## this runtime's width-rounding), any ColumnWidth in roughly (39.08, 39.25]
## rounds to an OOXML width of 40 - 39.1666... is the middle of that band.
$targetColumnWidth = 39.1666666667

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("P3").Value = "Handback file name: kahpvanb.oqy is different with handoff file name: 65f37892-d67e-47ee-a87e-c6f046ea05a5.a90d37d26714b6328d3e623ecd34be662b33ab07.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = $targetColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("P3").Value = "Handback file name: kahpvanb.oqy is different with handoff file name: 65f37892-d67e-47ee-a87e-c6f046ea05a5.a90d37d26714b6328d3e623ecd34be662b33ab07.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = $targetColumnWidth
